# Add a new manifest row describing the supplemental PDF file that was
# included alongside review-results.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "supplemental.pdf"
$ws.Range("B3").Value = "2021-09-07T12:04:41,040177-07:00"
$ws.Range("C3").Value = "Contains supplementary paragraphs and  tables referred in our paper"
$ws.Range("D3").Value = ".pdf"

# The "description" column no longer needs to be as wide once the very long
# description text is no longer the narrowest-fitting entry driving autosize.
$ws.Columns.Item(3).ColumnWidth = 62.67

# Leave the selection on the (still long) description cell of the first row,
# matching where the author's cursor ended up after the edit.
$ws.Range("C2").Select()
